$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 9.718953170468685
$ws.Range("C2").Value = -8.452308629668783
$ws.Range("D2").Value = -0.3126370826896341
$ws.Range("E2").Value = 1.834933389109758
$ws.Range("F2").Value = -1.605096597413572
$ws.Range("G2").Value = -1.666270775827773
$ws.Range("H2").Value = 0.3751851613284818
$ws.Range("I2").Value = -0.2074093757051592
$ws.Range("J2").Value = -0.1117343348958014
$ws.Range("K2").Value = -0.3721533718945242

# Row 3
$ws.Range("B3").Value = -10.70211146928018
$ws.Range("C3").Value = -2.562439922301026
$ws.Range("D3").Value = -0.4148694505016339
$ws.Range("E3").Value = -3.854899437024964
$ws.Range("F3").Value = -3.916073615439165
$ws.Range("G3").Value = -1.87461767828291
$ws.Range("H3").Value = -2.457212215316551
$ws.Range("I3").Value = -2.361537174507194
$ws.Range("J3").Value = -2.621956211505916
$ws.Range("K3").Value = -1.786501215347425

# Row 4
$ws.Range("B4").Value = -10.03159025315832
$ws.Range("C4").Value = -7.884019781358927
$ws.Range("D4").Value = -11.32404976788226
$ws.Range("E4").Value = -11.38522394629646
$ws.Range("F4").Value = -9.343768009140204
$ws.Range("G4").Value = -9.926362546173845
$ws.Range("H4").Value = -9.830687505364487
$ws.Range("I4").Value = -10.09110654236321
$ws.Range("J4").Value = -9.255651546204719
$ws.Range("K4").Value = -9.35573090610734

# Row 5
$ws.Range("B5").Value = 10.28724201877854
$ws.Range("C5").Value = 6.847212032255212
$ws.Range("D5").Value = 6.786037853841011
$ws.Range("E5").Value = 8.827493790997265
$ws.Range("F5").Value = 8.244899253963624
$ws.Range("G5").Value = 8.340574294772981
$ws.Range("H5").Value = 8.080155257774258
$ws.Range("I5").Value = 8.91561025393275
$ws.Range("J5").Value = 8.815530894030129
$ws.Range("K5").Value = 8.461744387950432

# Row 6
$ws.Range("B6").Value = -1.292459514723937
$ws.Range("C6").Value = -1.353633693138139
$ws.Range("D6").Value = 0.6878222440181159
$ws.Range("E6").Value = 0.1052277069844749
$ws.Range("F6").Value = 0.2009027477938327
$ws.Range("G6").Value = -0.0595162892048901
$ws.Range("H6").Value = 0.7759387069536011
$ws.Range("I6").Value = 0.6758593470509796
$ws.Range("J6").Value = 0.3220728409712834
$ws.Range("K6").Value = 0.4702502314037287

# Row 7
$ws.Range("B7").Value = -3.501204164937531
$ws.Range("C7").Value = -1.459748227781277
$ws.Range("D7").Value = -2.042342764814918
$ws.Range("E7").Value = -1.94666772400556
$ws.Range("F7").Value = -2.207086761004283
$ws.Range("G7").Value = -1.371631764845791
$ws.Range("H7").Value = -1.471711124748413
$ws.Range("I7").Value = -1.825497630828109
$ws.Range("J7").Value = -1.677320240395664
$ws.Range("K7").Value = -1.519874251391509

# Row 8
$ws.Range("B8").Value = 1.980281758742053
$ws.Range("C8").Value = 1.397687221708412
$ws.Range("D8").Value = 1.49336226251777
$ws.Range("E8").Value = 1.232943225519047
$ws.Range("F8").Value = 2.068398221677539
$ws.Range("G8").Value = 1.968318861774917
$ws.Range("H8").Value = 1.614532355695221
$ws.Range("I8").Value = 1.762709746127666
$ws.Range("J8").Value = 1.920155735131821
$ws.Range("K8").Value = 1.41949627860575

# Row 9
$ws.Range("B9").Value = 1.458861400122613
$ws.Range("C9").Value = 1.554536440931971
$ws.Range("D9").Value = 1.294117403933248
$ws.Range("E9").Value = 2.12957240009174
$ws.Range("F9").Value = 2.029493040189118
$ws.Range("G9").Value = 1.675706534109422
$ws.Range("H9").Value = 1.823883924541867
$ws.Range("I9").Value = 1.981329913546022
$ws.Range("J9").Value = 1.480670457019951
$ws.Range("K9").Value = 1.761028337555421

# Row 10
$ws.Range("B10").Value = -0.4869194962242832
$ws.Range("C10").Value = -0.747338533223006
$ws.Range("D10").Value = 0.0881164629354852
$ws.Range("E10").Value = -0.01196289696713632
$ws.Range("F10").Value = -0.3657494030468326
$ws.Range("G10").Value = -0.2175720126143872
$ws.Range("H10").Value = -0.06012602361023223
$ws.Range("I10").Value = -0.560785480136303
$ws.Range("J10").Value = -0.2804275996008339
$ws.Range("K10").Value = -0.4524024574458455

# Row 11
$ws.Range("B11").Value = -0.164743996189365
$ws.Range("C11").Value = 0.6707109999691262
$ws.Range("D11").Value = 0.5706316400665047
$ws.Range("E11").Value = 0.2168451339868084
$ws.Range("F11").Value = 0.3650225244192538
$ws.Range("G11").Value = 0.5224685134234088
$ws.Range("H11").Value = 0.02180905689733798
$ws.Range("I11").Value = 0.3021669374328071
$ws.Range("J11").Value = 0.1301920795877955
$ws.Range("K11").Value = 0.436054619334127

# Row 12
$ws.Range("B12").Value = 0.5750359591597685
$ws.Range("C12").Value = 0.4749565992571469
$ws.Range("D12").Value = 0.1211700931774507
$ws.Range("E12").Value = 0.269347483609896
$ws.Range("F12").Value = 0.426793472614051
$ws.Range("G12").Value = -0.07386598391201982
$ws.Range("H12").Value = 0.2064918966234494
$ws.Range("I12").Value = 0.0345170387784377
$ws.Range("J12").Value = 0.3403795785247692
$ws.Range("K12").Value = -0.271950374162517

# Row 13
$ws.Range("B13").Value = 0.7353756362558697
$ws.Range("C13").Value = 0.3815891301761735
$ws.Range("D13").Value = 0.5297665206086188
$ws.Range("E13").Value = 0.6872125096127738
$ws.Range("F13").Value = 0.186553053086703
$ws.Range("G13").Value = 0.4669109336221722
$ws.Range("H13").Value = 0.2949360757771605
$ws.Range("I13").Value = 0.600798615523492
$ws.Range("J13").Value = -0.01153133716379418
$ws.Range("K13").Value = 0.6768400480353174

# Row 14
$ws.Range("B14").Value = -0.4538658659823178
$ws.Range("C14").Value = -0.3056884755498724
$ws.Range("D14").Value = -0.1482424865457174
$ws.Range("E14").Value = -0.6489019430717882
$ws.Range("F14").Value = -0.3685440625363191
$ws.Range("G14").Value = -0.5405189203813308
$ws.Range("H14").Value = -0.2346563806349992
$ws.Range("I14").Value = -0.8469863333222853
$ws.Range("J14").Value = -0.1586149481231739
$ws.Range("K14").Value = -0.4353061035472806

# Row 15
$ws.Range("B15").Value = -0.2056091156472509
$ws.Range("C15").Value = -0.04816312664309591
$ws.Range("D15").Value = -0.5488225831691667
$ws.Range("E15").Value = -0.2684647026336975
$ws.Range("F15").Value = -0.4404395604787092
$ws.Range("G15").Value = -0.1345770207323777
$ws.Range("H15").Value = -0.7469069734196638
$ws.Range("I15").Value = -0.05853558822055238
$ws.Range("J15").Value = -0.3352267436446591
$ws.Range("K15").ClearContents()

# Row 16
$ws.Range("B16").Value = 0.3056233794366003
$ws.Range("C16").Value = -0.1950360770894705
$ws.Range("D16").Value = 0.08532180344599868
$ws.Range("E16").Value = -0.08665305439901295
$ws.Range("F16").Value = 0.2192094853473185
$ws.Range("G16").Value = -0.3931204673399676
$ws.Range("H16").Value = 0.2952509178591439
$ws.Range("I16").Value = 0.01855976243503714
$ws.Range("J16").ClearContents()

# Row 17
$ws.Range("B17").Value = -0.3432134675219158
$ws.Range("C17").Value = -0.06285558698644665
$ws.Range("D17").Value = -0.2348304448314583
$ws.Range("E17").Value = 0.0710320949148732
$ws.Range("F17").Value = -0.541297857772413
$ws.Range("G17").Value = 0.1470735274266985
$ws.Range("H17").Value = -0.1296176279974082
$ws.Range("I17").ClearContents()

# Row 18
$ws.Range("B18").Value = -0.2203015759906016
$ws.Range("C18").Value = -0.3922764338356133
$ws.Range("D18").Value = -0.0864138940892818
$ws.Range("E18").Value = -0.698743846776568
$ws.Range("F18").Value = -0.01037246157745647
$ws.Range("G18").Value = -0.2870636170015632
$ws.Range("H18").ClearContents()

# Row 19
$ws.Range("B19").Value = 0.1083830226904575
$ws.Range("C19").Value = 0.414245562436789
$ws.Range("D19").Value = -0.1980843902504972
$ws.Range("E19").Value = 0.4902869949486143
$ws.Range("F19").Value = 0.2135958395245076
$ws.Range("G19").ClearContents()

# Row 20
$ws.Range("B20").Value = 0.1338876819013198
$ws.Range("C20").Value = -0.4784422707859664
$ws.Range("D20").Value = 0.2099291144131452
$ws.Range("E20").Value = -0.06676204101096155
$ws.Range("F20").ClearContents()

# Row 21
$ws.Range("B21").Value = -0.3064674129409547
$ws.Range("C21").Value = 0.3819039722581568
$ws.Range("D21").Value = 0.1052128168340501
$ws.Range("E21").ClearContents()

# Row 22
$ws.Range("B22").Value = 0.07604143251182532
$ws.Range("C22").Value = -0.2006497229122814
$ws.Range("D22").ClearContents()

# Row 23
$ws.Range("B23").Value = 0.4116802297750048
$ws.Range("C23").ClearContents()

# Row 24
$ws.Range("B24").ClearContents()
